# A new weekly price record for "Alcachofa" (Madrigal / Primera) is inserted
# at row 42 of the sheet, pushing all the existing rows (old 42..131) down by
# one row (to 43..132). This matches the commit message
# "Fruta / hortaliza, semanal" (a new weekly row of fruit/vegetable data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row before the current row 42 - this shifts rows
# 42..131 down to 43..132 and keeps formatting (incl. the date style used
# in column D) carried over from the row above, just like Excel's UI does.
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new record's data.
$ws.Range("A42").Value = 5
$ws.Range("B42").Value = "Macroferia Regional de Talca"
$ws.Range("C42").Value = "Maule"
$ws.Range("D42").Value = 45152
$ws.Range("E42").Value = 7
$ws.Range("F42").Value = 100112013
$ws.Range("G42").Value = "Alcachofa"
$ws.Range("H42").Value = "Madrigal"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 600
$ws.Range("K42").Value = 13000
$ws.Range("L42").Value = 13000
$ws.Range("M42").Value = 13000
$ws.Range("N42").Value = "`$/caja 40 unidades"
$ws.Range("O42").Value = "Provincia del Elquí"
$ws.Range("P42").Value = 325
$ws.Range("Q42").Value = 40
$ws.Range("R42").Value = "Hortaliza"

Write-Output "inserted new row 42"
